$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "27.742.42"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.63%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.636.60"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -0.63%  "
$ws.Range("E4").Value = "  -0.03%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "212.55"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.06%  "
$ws.Range("E6").Value = "  -1.89%  "
$ws.Range("E7").Value = "  -0.06%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "23.29"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -0.96%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.261"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +1.96%  "
$ws.Range("E10").Value = "  -0.13%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0888"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.08%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.869.52"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.64%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.639.15"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -0.49%  "
$ws.Range("E14").Value = "  +0.25%  "
$ws.Range("E15").Value = "  -3.89%  "
$ws.Range("E16").Value = "  +0.54%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "27.733.09"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.69%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "229.92"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.68%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "7.69"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +1.52%  "
$ws.Range("E20").Value = "  -0.10%  "
$ws.Range("E21").Value = "  +0.00%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.30"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.65%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "10.22"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +4.61%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.08"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +4.38%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "151.24"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +1.37%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "6.93"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.99%  "
$ws.Range("E27").Value = "  -0.73%  "
$ws.Range("E28").Value = "  -0.09%  "
$ws.Range("E29").Value = "  -0.02%  "
$ws.Range("E30").Value = "  +0.25%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.0485"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.12%  "
$ws.Range("E32").Value = "  +0.06%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.469.52"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +2.71%  "
$ws.Range("E34").Value = "  -2.55%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.56"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -2.95%  "
$ws.Range("E36").Value = "  -0.74%  "
$ws.Range("E37").Value = "  -0.55%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.879"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -0.46%  "
$ws.Range("E39").Value = "  +0.02%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.902"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +10.17%  "
$ws.Range("B41").Value = "PaxDollar"
$ws.Range("C41").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.00"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.04%  "
$ws.Range("B42").Value = "Aave"
$ws.Range("C42").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "68.61"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +5.48%  "
$ws.Range("E43").Value = "  -1.85%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "5.59"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +1.23%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.47"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.01%  "
$ws.Range("E46").Value = "  -0.98%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.779.28"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.62%  "
$ws.Range("E48").Value = "  +1.45%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "87.09"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -1.14%  "
$ws.Range("E50").Value = "  -1.50%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0991"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.65%  "
